{"js": "// The document contains a single table of 20 rows x 5 columns, where\n// every cell holds one arithmetic expression (e.g. \"90-8=\"). The edit\n// replaces the expression text in every cell with a new expression,\n// while leaving the cell/paragraph/run formatting untouched.\n//\n// Writing `Table.values` rewrites only the text of the existing runs\n// (it does not touch paragraph/run formatting such as font or size),\n// so it's the safest way to do a wholesale content swap cell-by-cell\n// in document (row-major) order.\n\nconst newValues = [\n  [\"42+22=\", \"51-51=\", \"91+7=\", \"33+44=\", \"36+44=\"],\n  [\"60+14=\", \"46+21=\", \"58-24=\", \"76-39=\", \"25+61=\"],\n  [\"40+38=\", \"68-30=\", \"94-64=\", \"82-34=\", \"43+49=\"],\n  [\"95+1=\", \"14-4=\", \"61+24=\", \"90-21=\", \"39-22=\"],\n  [\"6+93=\", \"19+3=\", \"53+7=\", \"6+7=\", \"12+83=\"],\n  [\"30+2=\", \"57-41=\", \"8+55=\", \"96-60=\", \"16+67=\"],\n  [\"52+31=\", \"98-72=\", \"28-16=\", \"64+16=\", \"26-1=\"],\n  [\"73+1=\", \"65+31=\", \"37-18=\", \"80-56=\", \"5+43=\"],\n  [\"81+5=\", \"59+32=\", \"2+25=\", \"6-4=\", \"70-62=\"],\n  [\"8+3=\", \"55-33=\", \"29+4=\", \"46+4=\", \"35+8=\"],\n  [\"53-18=\", \"99-0=\", \"20+17=\", \"58+26=\", \"78-20=\"],\n  [\"14+61=\", \"39-37=\", \"74-49=\", \"21+19=\", \"86-15=\"],\n  [\"25+26=\", \"0+16=\", \"96-69=\", \"87+0=\", \"41+54=\"],\n  [\"56+20=\", \"6+33=\", \"26+65=\", \"82-78=\", \"35-11=\"],\n  [\"49+40=\", \"49-48=\", \"28-15=\", \"83+3=\", \"46+26=\"],\n  [\"98-85=\", \"28+21=\", \"11-2=\", \"58+2=\", \"44-35=\"],\n  [\"56-10=\", \"24+40=\", \"26+47=\", \"29+21=\", \"14-4=\"],\n  [\"70-6=\", \"37-30=\", \"80-71=\", \"55-52=\", \"96-42=\"],\n  [\"12+57=\", \"42-20=\", \"94-54=\", \"87-60=\", \"63-45=\"],\n  [\"39+27=\", \"29+64=\", \"83-75=\", \"1+21=\", \"43+14=\"],\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nif (table.rowCount !== newValues.length) {\n  throw new Error(\n    `Expected ${newValues.length} rows, found ${table.rowCount}`\n  );\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# The document contains a single table of 20 rows x 5 columns, where\n# every cell holds one arithmetic expression (e.g. \"90-8=\"). The edit\n# replaces the expression text in every cell with a new expression, in\n# row-major document order, while leaving the cell/paragraph/run\n# formatting untouched (only Range.Text is reassigned).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    '42+22=', '51-51=', '91+7=', '33+44=', '36+44=',\n    '60+14=', '46+21=', '58-24=', '76-39=', '25+61=',\n    '40+38=', '68-30=', '94-64=', '82-34=', '43+49=',\n    '95+1=', '14-4=', '61+24=', '90-21=', '39-22=',\n    '6+93=', '19+3=', '53+7=', '6+7=', '12+83=',\n    '30+2=', '57-41=', '8+55=', '96-60=', '16+67=',\n    '52+31=', '98-72=', '28-16=', '64+16=', '26-1=',\n    '73+1=', '65+31=', '37-18=', '80-56=', '5+43=',\n    '81+5=', '59+32=', '2+25=', '6-4=', '70-62=',\n    '8+3=', '55-33=', '29+4=', '46+4=', '35+8=',\n    '53-18=', '99-0=', '20+17=', '58+26=', '78-20=',\n    '14+61=', '39-37=', '74-49=', '21+19=', '86-15=',\n    '25+26=', '0+16=', '96-69=', '87+0=', '41+54=',\n    '56+20=', '6+33=', '26+65=', '82-78=', '35-11=',\n    '49+40=', '49-48=', '28-15=', '83+3=', '46+26=',\n    '98-85=', '28+21=', '11-2=', '58+2=', '44-35=',\n    '56-10=', '24+40=', '26+47=', '29+21=', '14-4=',\n    '70-6=', '37-30=', '80-71=', '55-52=', '96-42=',\n    '12+57=', '42-20=', '94-54=', '87-60=', '63-45=',\n    '39+27=', '29+64=', '83-75=', '1+21=', '43+14='\n)\n\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\nif (($rows * $cols) -ne $newValues.Count) {\n    throw \"Expected $($newValues.Count) cells, found $($rows * $cols)\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n"}
